# Re-pull data / push all data / mean calculation
# Update the dSF column (F) values to reflect freshly pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    4  = -1
    5  = 1
    6  = 4
    8  = 1
    9  = -5
    10 = 1
    11 = 2
    13 = 2
    14 = -2
    15 = -1
    16 = 0
    18 = 4
    19 = 5
    21 = -2
    22 = -2
    23 = 5
    25 = -2
    29 = -1
    31 = -4
    32 = 2
    36 = -6
    37 = -1
    38 = -1
    39 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
